# Weekly refresh: insert the newest Ciboulette price record for
# "Vega Modelo de Temuco" ahead of the existing history (row 233),
# pushing the previously recorded rows (233:250) down by one row
# (234:251), exactly like a new week's row being prepended to this
# market/category block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 233; Excel copies formatting (incl. the date
# style on column D) down from the row above, matching the workbook's
# existing pattern.
$ws.Rows("233:233").Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(233, 1).Value = 10
$ws.Cells.Item(233, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(233, 3).Value = "La Araucanía"
$ws.Cells.Item(233, 4).Value = 44714
$ws.Cells.Item(233, 5).Value = 9
$ws.Cells.Item(233, 6).Value = 100112039
$ws.Cells.Item(233, 7).Value = "Ciboulette"
$ws.Cells.Item(233, 8).Value = "Sin especificar"
$ws.Cells.Item(233, 9).Value = "Primera"
$ws.Cells.Item(233, 10).Value = 65
$ws.Cells.Item(233, 11).Value = 6000
$ws.Cells.Item(233, 12).Value = 6000
$ws.Cells.Item(233, 13).Value = 6000
$ws.Cells.Item(233, 14).Value = "`$/docena de atados"
$ws.Cells.Item(233, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(233, 16).Value = 2000
$ws.Cells.Item(233, 17).Value = 3
$ws.Cells.Item(233, 18).Value = "Hortaliza"
